# Applies the marksheet re-grading update to the quiz result sheet.
# - Updates the summary table (rows 10-12) with real counts/marks instead of the
#   "Absent" placeholder values.
# - Fills in the student's answers for the first question block (column A,
#   rows 16-40), coloring each one green ("correctStyle") when it matches the
#   correct answer in column B, or red ("incorrectStyle") when it doesn't;
#   unattempted questions stay blank with the plain "normalStyle".
# - Fills the first three rows of the second question block (columns D/E) the
#   same way, and clears out the remainder of that block as well as the whole
#   third block (columns G/H), which are no longer part of the report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Copy-StyleTo($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats) | Out-Null
}

# ---- Summary table (rows 10-12) ----------------------------------------

# Make the row-label cells use the same bold "mtitleStyle" formatting as the
# other header cells (A15, B15, ...).
Copy-StyleTo "A15" "A10:A12"

$ws.Range("A10").Value = "No."
$ws.Range("B10").Value = 11
$ws.Range("C10").Value = 6
$ws.Range("D10").Value = 11
$ws.Range("E10").Value = 28

$ws.Range("A11").Value = "Marking"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 0

$ws.Range("A12").Value = "Total"
$ws.Range("B12").Value = 44
$ws.Range("C12").Value = -6
$ws.Range("E12").Value = "38/112"

# ---- Question block 3 (G/H) is removed entirely -------------------------

$ws.Range("G15:H40").Clear() | Out-Null

# ---- Question block 2 (D/E): keep only the first 3 answered rows --------

$ws.Range("D16").Value = "Option A"
Copy-StyleTo "B10" "D16"   # matches E16 ("Option A") -> correct/green

$ws.Range("D17").Value = "Option B"
Copy-StyleTo "C10" "D17"   # does not match E17 ("Option C") -> incorrect/red

$ws.Range("D18").Value = "Option D"
Copy-StyleTo "B10" "D18"   # matches E18 ("Option D") -> correct/green

$ws.Range("D19:E40").Clear() | Out-Null

# ---- Question block 1 (A/B): fill student answers in column A -----------

$correctRows    = @(17, 18, 27, 28, 32, 33, 37, 38, 40)
$incorrectRows  = @(16, 19, 22, 26, 29)

$answers = @{
    16 = "Option D"
    17 = "Option D"
    18 = "Option B"
    19 = "Option D"
    22 = "Option B"
    26 = "Option A"
    27 = "Option A"
    28 = "Option D"
    29 = "Option A"
    32 = "Option C"
    33 = "Option D"
    37 = "Option A"
    38 = "Option A"
    40 = "Option D"
}

foreach ($r in $answers.Keys) {
    $ws.Range("A$r").Value = $answers[$r]
}

foreach ($r in $correctRows) {
    Copy-StyleTo "B10" "A$r"
}

foreach ($r in $incorrectRows) {
    Copy-StyleTo "C10" "A$r"
}
